$d = $word.ActiveDocument

$d.Content.Find.Execute("529×5=2645", $true, $false, $false, $false, $false, $true, 1, $false, "683×2=1366", 2)
$d.Content.Find.Execute("422×7=2954", $true, $false, $false, $false, $false, $true, 1, $false, "669×9=6021", 2)
$d.Content.Find.Execute("841×5=4205", $true, $false, $false, $false, $false, $true, 1, $false, "324×4=1296", 2)
$d.Content.Find.Execute("256×3=768", $true, $false, $false, $false, $false, $true, 1, $false, "911×3=2733", 2)
$d.Content.Find.Execute("542×3=1626", $true, $false, $false, $false, $false, $true, 1, $false, "423×4=1692", 2)
$d.Content.Find.Execute("230×9=2070", $true, $false, $false, $false, $false, $true, 1, $false, "564×8=4512", 2)
$d.Content.Find.Execute("689×4=2756", $true, $false, $false, $false, $false, $true, 1, $false, "269×3=807", 2)
$d.Content.Find.Execute("243×3=729", $true, $false, $false, $false, $false, $true, 1, $false, "389×3=1167", 2)
$d.Content.Find.Execute("701×5=3505", $true, $false, $false, $false, $false, $true, 1, $false, "968×4=3872", 2)
$d.Content.Find.Execute("779×8=6232", $true, $false, $false, $false, $false, $true, 1, $false, "825×6=4950", 2)
$d.Content.Find.Execute("527×5=2635", $true, $false, $false, $false, $false, $true, 1, $false, "302×6=1812", 2)
$d.Content.Find.Execute("515×8=4120", $true, $false, $false, $false, $false, $true, 1, $false, "441×2=882", 2)
$d.Content.Find.Execute("866×2=1732", $true, $false, $false, $false, $false, $true, 1, $false, "501×7=3507", 2)
$d.Content.Find.Execute("378×5=1890", $true, $false, $false, $false, $false, $true, 1, $false, "625×2=1250", 2)
$d.Content.Find.Execute("774×4=3096", $true, $false, $false, $false, $false, $true, 1, $false, "475×5=2375", 2)
$d.Content.Find.Execute("333×6=1998", $true, $false, $false, $false, $false, $true, 1, $false, "844×9=7596", 2)
$d.Content.Find.Execute("702×6=4212", $true, $false, $false, $false, $false, $true, 1, $false, "878×7=6146", 2)
$d.Content.Find.Execute("476×4=1904", $true, $false, $false, $false, $false, $true, 1, $false, "321×8=2568", 2)
$d.Content.Find.Execute("755×5=3775", $true, $false, $false, $false, $false, $true, 1, $false, "137×5=685", 2)
$d.Content.Find.Execute("247×5=1235", $true, $false, $false, $false, $false, $true, 1, $false, "415×7=2905", 2)
$d.Content.Find.Execute("259×4=1036", $true, $false, $false, $false, $false, $true, 1, $false, "742×4=2968", 2)
$d.Content.Find.Execute("952×7=6664", $true, $false, $false, $false, $false, $true, 1, $false, "675×4=2700", 2)
$d.Content.Find.Execute("200×2=400", $true, $false, $false, $false, $false, $true, 1, $false, "498×9=4482", 2)
$d.Content.Find.Execute("548×8=4384", $true, $false, $false, $false, $false, $true, 1, $false, "115×5=575", 2)
$d.Content.Find.Execute("439×5=2195", $true, $false, $false, $false, $false, $true, 1, $false, "656×5=3280", 2)
